$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data right
$ws.Columns("A:A").Insert()

# New column A values
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "Cases"

# Column width (only the new column A needs a width set; B:E already
# retain their original widths after the insert shifted them right).
# (target char-width is 8.81640625; feeding that value back through the
# ColumnWidth setter's internal pixel rounding lands closest at this input)
$ws.Columns("A:A").ColumnWidth = 7.983072916666667

# Update selection
$ws.Range("B8").Select()
